# Applies the "add end-to-end test and faker data class" edit:
#  - register sheet (sheet1): test data row for username/email refreshed
#    from "iti_5" / "iti_5@iti.com" to "iti_7" / "iti_7@iti.com"; this is
#    now the sheet the user last looked at (tabSelected + selection B2)
#  - categoryNames sheet (sheet4): removed the no-longer-needed
#    Leggings / Jumpers / Growsuits / Dresses rows
#  - postsNames sheet (sheet6): removed the no-longer-needed
#    "growsuit", "sleeping bag" and "pain relief" post-name rows, leaving
#    the "organic cotton" post name as the new row 3
#  - productsNames (2) sheet (sheet8): no longer the active tab

$wb = $excel.ActiveWorkbook

# --- register sheet: refresh the faker test data row -----------------
$wsRegister = $wb.Worksheets.Item("register")
$wsRegister.Range("A2").Value = "iti_7"
$wsRegister.Range("B2").Value = "iti_7@iti.com"

# --- categoryNames sheet: drop the unused category rows ---------------
$wsCategory = $wb.Worksheets.Item("categoryNames")
$wsCategory.Rows("3:6").Delete()
$wsCategory.Rows(3).Select()

# --- postsNames sheet: drop the unused post-name rows ------------------
$wsPosts = $wb.Worksheets.Item("postsNames")
$wsPosts.Rows("3:5").Delete()
$wsPosts.Rows(3).Select()

# --- make "register" the active/selected sheet+cell -------------------
$wsRegister.Select()
$wsRegister.Range("B2").Select()
